$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (not General/number) formatting on Price column cells we are
# about to rewrite, so numeric-looking strings (e.g. "1.00", "0.0000242",
# "2.30") are preserved verbatim as text instead of being normalised to a
# number (matching the original inline-string cell type).
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"

# Apply the updated values scraped by the GitHub Actions job.
$ws.Range('D2').Value = '67.072.52'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '3.123.00'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '579.44'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').Value = '174.85'
$ws.Range('E6').Value = '  +0.39%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.122.56'
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').Value = '0.517'
$ws.Range('E9').Value = '  -1.07%  '
$ws.Range('D10').Value = '6.41'
$ws.Range('E10').Value = '  +0.14%  '
$ws.Range('D11').Value = '0.153'
$ws.Range('E11').Value = '  -1.83%  '
$ws.Range('E12').Value = '  -1.17%  '
$ws.Range('D13').Value = '0.0000242'
$ws.Range('E13').Value = '  -3.16%  '
$ws.Range('D14').Value = '36.23'
$ws.Range('E14').Value = '  -2.81%  '
$ws.Range('E15').Value = '  -0.45%  '
$ws.Range('D16').Value = '3.641.10'
$ws.Range('E16').Value = '  +0.00%  '
$ws.Range('D17').Value = '67.023.35'
$ws.Range('E17').Value = '  -0.10%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').Value = '17.11'
$ws.Range('E18').Value = '  +4.04%  '
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').Value = '7.02'
$ws.Range('E19').Value = '  -1.56%  '
$ws.Range('D20').Value = '3.112.12'
$ws.Range('E20').Value = '  -0.38%  '
$ws.Range('D21').Value = '486.87'
$ws.Range('E21').Value = '  -0.80%  '
$ws.Range('D22').Value = '7.88'
$ws.Range('E22').Value = '  -0.98%  '
$ws.Range('D23').Value = '0.695'
$ws.Range('E23').Value = '  -1.82%  '
$ws.Range('D24').Value = '83.89'
$ws.Range('E24').Value = '  -0.37%  '
$ws.Range('D25').Value = '12.84'
$ws.Range('E25').Value = '  -2.93%  '
$ws.Range('D26').Value = '2.25'
$ws.Range('E26').Value = '  -1.99%  '
$ws.Range('D27').Value = '10.23'
$ws.Range('E27').Value = '  -1.37%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('D29').Value = '8.04'
$ws.Range('E29').Value = '  +1.69%  '
$ws.Range('D30').Value = '2.30'
$ws.Range('E30').Value = '  -2.44%  '
$ws.Range('D31').Value = '2.61'
$ws.Range('E31').Value = '  -2.62%  '
$ws.Range('D32').Value = '28.14'
$ws.Range('E32').Value = '  -1.72%  '
$ws.Range('D33').Value = '0.113'
$ws.Range('E33').Value = '  -1.78%  '
$ws.Range('D34').Value = '0.0₃0946'
$ws.Range('E34').Value = '  -0.60%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('E36').Value = '  +1.63%  '
$ws.Range('E37').Value = '  -4.45%  '
$ws.Range('E38').Value = '  -3.05%  '
$ws.Range('D39').Value = '0.313'
$ws.Range('E39').Value = '  +0.52%  '
$ws.Range('D40').Value = '49.30'
$ws.Range('E40').Value = '  -1.67%  '
$ws.Range('E41').Value = '  -0.44%  '
$ws.Range('E42').Value = '  -3.45%  '
$ws.Range('D43').Value = '8.34'
$ws.Range('E43').Value = '  -2.29%  '
$ws.Range('D44').Value = '2.68'
$ws.Range('E44').Value = '  +2.39%  '
$ws.Range('D45').Value = '2.812.52'
$ws.Range('E45').Value = '  -0.46%  '
$ws.Range('D46').Value = '0.0349'
$ws.Range('E46').Value = '  -1.26%  '
$ws.Range('D47').Value = '373.55'
$ws.Range('E47').Value = '  -2.86%  '
$ws.Range('D48').Value = '134.78'
$ws.Range('E48').Value = '  -0.56%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').Value = '24.73'
$ws.Range('E50').Value = '  -1.09%  '
$ws.Range('E51').Value = '  +0.50%  '
